$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell contents for rows 10-23 to match the target layout
$ws.Range("A10").Value = "Objetivos:"
$ws.Range("B10").Value = "11079086 - Herlandí de Souza Andrade"
$ws.Range("C10").Value = "11079086 - Herlandí de Souza Andrade"

$ws.Range("A11").Value = "Objectives:"
$ws.Range("B11").Value = "Present the fundamentals of General Accounting, so that students can understand the basic concepts and principles of Accounting. Train in the production and analysis of documents and records of accounting transactions. Understand costing systems and their application. Empower to design and implement costing systems. Empower to analyze costs in relation to an organization's market and production strategy."
$ws.Range("C11").Value = "Present the fundamentals of General Accounting, so that students can understand the basic concepts and principles of Accounting. Train in the production and analysis of documents and records of accounting transactions. Understand costing systems and their application. Empower to design and implement costing systems. Empower to analyze costs in relation to an organization's market and production strategy."

$ws.Range("A12").Value = "Docentes responsáveis:"
$ws.Range("B12").ClearContents() | Out-Null
$ws.Range("C12").ClearContents() | Out-Null

$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = "1 - Fundamentals of Accounting. 2 - Costing Systems. 3 - Costing Systems Project. 4 - Cost Analysis."
$ws.Range("C14").Value = "1 - Fundamentals of Accounting. 2 - Costing Systems. 3 - Costing Systems Project. 4 - Cost Analysis."

$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "01/01/2021"
$ws.Range("C15").Value = "01/01/2021"

$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = "1. ACCOUNTING FUNDAMENTALSAccounting concept. Production Accounting and Engineering. Balance sheet: assets, liabilities, equity. Basic accounting procedures. Changes in equity; expense, revenue. Accrual and cash regimes. Deferred income and expenses. Accounting Facts. Goods operations. Inventory. Property, Plant and Equipment and Amortization. Statement of origins and applications. Balance sheet and results analysis.2. COSTING SYSTEMSDirect and indirect costs, fixed and variable. Accumulation of costs, classification, cost factors. Absorption costing. Orderly, continuous, joint production. Direct Costing: contribution margin. ABC costing.3. COSTING SYSTEMS DESIGNProducts and Departments. Production and costing systems. The inflation problem.4. COST ANALYSISCost-Volume-Profit. Marginal contribution. Analysis of variations. Global productivity equation. Operational leverage. IRR and profitability."
$ws.Range("C16").Value = "1. ACCOUNTING FUNDAMENTALSAccounting concept. Production Accounting and Engineering. Balance sheet: assets, liabilities, equity. Basic accounting procedures. Changes in equity; expense, revenue. Accrual and cash regimes. Deferred income and expenses. Accounting Facts. Goods operations. Inventory. Property, Plant and Equipment and Amortization. Statement of origins and applications. Balance sheet and results analysis.2. COSTING SYSTEMSDirect and indirect costs, fixed and variable. Accumulation of costs, classification, cost factors. Absorption costing. Orderly, continuous, joint production. Direct Costing: contribution margin. ABC costing.3. COSTING SYSTEMS DESIGNProducts and Departments. Production and costing systems. The inflation problem.4. COST ANALYSISCost-Volume-Profit. Marginal contribution. Analysis of variations. Global productivity equation. Operational leverage. IRR and profitability."

$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17").ClearContents() | Out-Null
$ws.Range("C17").ClearContents() | Out-Null

$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "11079086 - Herlandí de Souza Andrade"
$ws.Range("C18").Value = "11079086 - Herlandí de Souza Andrade"

$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras."
$ws.Range("C19").Value = "Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras."

$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "Média Aritmética dos Projetos, Trabalhos e Exercícios realizados no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude) desenvolvidas."
$ws.Range("C20").Value = "Média Aritmética dos Projetos, Trabalhos e Exercícios realizados no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude) desenvolvidas."

$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação."
$ws.Range("C21").Value = "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação."

$ws.Range("A22").Value = "Requisitos:"
$ws.Range("B22").ClearContents() | Out-Null
$ws.Range("C22").ClearContents() | Out-Null

$ws.Range("A23").ClearContents() | Out-Null
$ws.Range("B23").Value = "LOQ4204 -  Economia Geral  (Requisito fraco)`n"
$ws.Range("C23").Value = "LOQ4204 -  Economia Geral  (Requisito fraco)`n"

# Remove the now-obsolete last row (content moved up into row 23)
$ws.Rows("24:24").Delete()

# Fix up row heights to match the target formatting
$ws.Rows(13).RowHeight = 60
$ws.Rows(15).RowHeight = 120
$ws.Rows(17).AutoFit()
$ws.Rows(18).RowHeight = 60
$ws.Rows(21).RowHeight = 120
$ws.Rows(22).AutoFit()
$ws.Rows(23).RowHeight = 30
